# "Split progression between inf and lowinf"
# Adds a new parameter row (progression_prop_infectious) right after
# progression_rate_age15, and tightens a few of the uniform-distribution
# upper bounds (distri_param2) used for sampling.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constant")
$ws.Activate()

# raw_transmission_rate: narrow the uniform prior (distri_param1/2)
$ws.Cells.Item(2, 4).Value() = 0.1
$ws.Cells.Item(2, 5).Value() = 10

# Insert a new row right after progression_rate_age15 (row 13) for the
# new parameter, inheriting the formatting of the row above it.
$ws.Rows.Item(14).Insert()
$ws.Cells.Item(14, 1).Value() = "progression_prop_infectious"
$ws.Cells.Item(14, 2).Value() = 0.5

# The rows below have all shifted down by one because of the insert.
# clinical_progression_rate (was row 19, now row 20): tighten distri_param2
$ws.Cells.Item(20, 5).Value() = 5

# infectiousness_gain_rate (was row 21, now row 22): tighten distri_param2
$ws.Cells.Item(22, 5).Value() = 5

# recent_detection_rate (was row 26, now row 27): tighten distri_param2
$ws.Cells.Item(27, 5).Value() = 5

# Update the saved selection/scroll position to match where the author
# was working when they saved.
[void]$ws.Range("E28").Select()
